# Daily attendance processing - re-sort the "Recorded By" (column G) values
# so the comma-separated list of recorders is in ordinal (case-sensitive)
# alphabetical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$comparer = [System.StringComparer]::Ordinal

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $orig = $cell.Value2

    if ($orig -ne $null -and $orig -ne "") {
        $parts = $orig -split ", "

        if ($parts.Count -gt 1) {
            $list = New-Object System.Collections.Generic.List[string]
            foreach ($p in $parts) { [void]$list.Add($p) }
            $list.Sort($comparer)
            $newVal = [string]::Join(", ", $list)

            if ($newVal -ne $orig) {
                $cell.Value2 = $newVal
            }
        }
    }
}
